$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.821.20'
$ws.Range("E2").Value = '  -1.43%  '

$ws.Range("D3").Value = '2.596.99'
$ws.Range("E3").Value = '  -2.05%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '552.86'
$ws.Range("E5").Value = '  +2.73%  '

$ws.Range("D6").Value = '143.22'
$ws.Range("E6").Value = '  -2.22%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  +5.01%  '

$ws.Range("D9").Value = '6.78'
$ws.Range("E9").Value = '  +0.45%  '

$ws.Range("E10").Value = '  -1.97%  '

$ws.Range("D12").Value = '0.335'
$ws.Range("E12").Value = '  -0.87%  '

$ws.Range("D13").Value = '3.055.83'
$ws.Range("E13").Value = '  -1.89%  '

$ws.Range("D14").Value = '58.776.27'
$ws.Range("E14").Value = '  -1.33%  '

$ws.Range("E15").Value = '  -2.14%  '

$ws.Range("D16").Value = '2.601.81'
$ws.Range("E16").Value = '  -2.12%  '

$ws.Range("E17").Value = '  -2.36%  '

$ws.Range("E18").Value = '  +1.13%  '

$ws.Range("D19").Value = '337.64'
$ws.Range("E19").Value = '  -0.75%  '

$ws.Range("D20").Value = '10.06'
$ws.Range("E20").Value = '  -2.66%  '

$ws.Range("E21").Value = '  -1.05%  '

$ws.Range("D23").Value = '66.81'
$ws.Range("E23").Value = '  +0.21%  '

$ws.Range("D24").Value = '0.428'
$ws.Range("E24").Value = '  +2.42%  '

$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("E26").Value = '  -3.55%  '

$ws.Range("E27").Value = '  -2.36%  '

$ws.Range("E28").Value = '  +0.58%  '

$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("E31").Value = '  +2.20%  '

$ws.Range("D32").Value = '154.51'
$ws.Range("E32").Value = '  +2.37%  '

$ws.Range("D33").Value = '18.95'
$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("E34").Value = '  -2.16%  '

$ws.Range("E35").Value = '  +6.21%  '

$ws.Range("D36").Value = '1.13'
$ws.Range("E36").Value = '  -1.08%  '

$ws.Range("D37").Value = '36.88'
$ws.Range("E37").Value = '  -1.17%  '

$ws.Range("E38").Value = '  +1.01%  '

$ws.Range("D39").Value = '0.837'
$ws.Range("E39").Value = '  -0.21%  '

$ws.Range("E40").Value = '  -0.41%  '

$ws.Range("D41").Value = '283.73'
$ws.Range("E41").Value = '  -1.22%  '

$ws.Range("D42").Value = '0.997'
$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("E43").Value = '  -1.72%  '

$ws.Range("D44").Value = '0.0955'
$ws.Range("E44").Value = '  +0.91%  '

$ws.Range("D45").Value = '10.63'
$ws.Range("E45").Value = '  -1.03%  '

$ws.Range("D46").Value = '0.0532'
$ws.Range("E46").Value = '  -1.49%  '

$ws.Range("E47").Value = '  -0.18%  '

$ws.Range("D48").Value = '1.945.82'
$ws.Range("E48").Value = '  -1.16%  '

$ws.Range("D49").Value = '118.53'
$ws.Range("E49").Value = '  +6.68%  '

$ws.Range("D50").Value = '17.90'
$ws.Range("E50").Value = '  -2.59%  '

$ws.Range("E51").Value = '  -5.83%  '
